# Apply crypto data refresh as described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.193.84'
$ws.Range('E2').Value = '  -3.26%  '
$ws.Range('D3').Value = '1.914.41'
$ws.Range('E3').Value = '  -4.12%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.96%  '
$ws.Range('D5').Value = "'327.78"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').Value = "'0.4681"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -6.06%  '
$ws.Range('D8').Value = "'0.4021"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.32%  '
$ws.Range('D9').Value = "'53.10"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.98%  '
$ws.Range('D10').Value = "'0.08417"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.22%  '
$ws.Range('D11').Value = "'1.044"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.88%  '
$ws.Range('D12').Value = "'22.09"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.93%  '
$ws.Range('D13').Value = '1.926.95'
$ws.Range('E13').Value = '  -3.70%  '
$ws.Range('D14').Value = "'7.434"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.55%  '
$ws.Range('D15').Value = "'6.064"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.86%  '
$ws.Range('D16').Value = "'1.004"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('D17').Value = "'89.63"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.27%  '
$ws.Range('D18').Value = "'0.00001067"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.87%  '
$ws.Range('D19').Value = "'0.06610"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').Value = "'17.99"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -8.38%  '
$ws.Range('D21').Value = "'0.9999"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('D22').Value = "'5.731"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('D23').Value = '28.219.00'
$ws.Range('E23').Value = '  -3.25%  '
$ws.Range('D24').Value = "'11.23"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.55%  '
$ws.Range('D25').Value = "'2.304"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('D26').Value = '2.166.32'
$ws.Range('E26').Value = '  -2.99%  '
$ws.Range('D27').Value = "'153.09"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.60%  '
$ws.Range('D28').Value = "'19.96"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.31%  '
$ws.Range('D29').Value = "'5.746"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.10%  '
$ws.Range('D30').Value = "'2.115"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.62%  '
$ws.Range('D31').Value = "'123.35"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.29%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = "'0.09628"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.49%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'0.9723"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.37%  '
$ws.Range('D34').Value = "'1.444"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.83%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = "'3.643"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = "'5.532"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.18%  '
$ws.Range('D37').Value = "'8.806"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.94%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'1.266"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.10%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.02294"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.49%  '
$ws.Range('D40').Value = "'0.06158"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.78%  '
$ws.Range('D41').Value = "'0.6142"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.52%  '
$ws.Range('D42').Value = "'11.01"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.38%  '
$ws.Range('D43').Value = "'1.000"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('D44').Value = "'0.1903"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.45%  '
$ws.Range('D45').Value = "'1.303"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.78%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.5843"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.15%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'12.81"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.17%  '
$ws.Range('D48').Value = "'2.020"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.82%  '
$ws.Range('D49').Value = "'3.427"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.90%  '
$ws.Range('D50').Value = "'0.06894"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('D51').Value = "'109.92"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.94%  '
